$d = $word.ActiveDocument

# 1) Fix "GPOI" typo -> "GPIO" in the main body text.
$d.Content.Find.Execute("GPOI", $true, $false, $false, $false, $false,
                         $true, 1, $false, "GPIO", 2)

# 2) Bump the cached footer date field result from 2021-01-25 to 2022-05-03.
# wdHeaderFooterPrimary = 1
$footerRange = $d.Sections(1).Footers(1).Range
$footerRange.Find.Execute("2021-01-25", $true, $false, $false, $false, $false,
                           $true, 1, $false, "2022-05-03", 2)
